$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3877.0908
$ws.Range("I32").Value = 3831.4614
$ws.Range("J32").Value = 3943
$ws.Range("K32").Value = 3831.4614
$ws.Range("L32").Value = 3943
$ws.Range("M32").Value = -3505.4614
$ws.Range("N32").Value = -4595

$ws.Range("H94").Value = 55700812
$ws.Range("I94").Value = 100012456
$ws.Range("J94").Value = 311251.5
$ws.Range("K94").Value = 100012456
$ws.Range("L94").Value = 311251.5
$ws.Range("M94").Value = -100012005
$ws.Range("N94").Value = -312153.5

$ws.Range("H113").Value = 10112.429
$ws.Range("I113").Value = 10739.833
$ws.Range("K113").Value = 10739.833
$ws.Range("M113").Value = -7485.833000000001

$ws.Range("H116").Value = 1238361.8
$ws.Range("I116").Value = 1589608.2
$ws.Range("K116").Value = 1589608.2
$ws.Range("M116").Value = -1586166.2

$ws.Range("H132").Value = 4125.2856
$ws.Range("I132").Value = 3907.7036
$ws.Range("K132").Value = 11723.1108
$ws.Range("M132").Value = -9193.110799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 763.3333
$ws.Range("I4").Value = 747.5
$ws.Range("K4").Value = 747.5
$ws.Range("M4").Value = -631.5

$ws.Range("H31").Value = 3782.4
$ws.Range("I31").Value = 3782.4
$ws.Range("K31").Value = 3782.4
$ws.Range("M31").Value = -3488.4

$ws.Range("H32").Value = 6140.525
$ws.Range("I32").Value = 6279.5
$ws.Range("K32").Value = 6279.5
$ws.Range("M32").Value = -5992.5

$ws.Range("H61").Value = 3643.5854
$ws.Range("I61").Value = 3003.4138
$ws.Range("J61").Value = 5190.6665
$ws.Range("K61").Value = 3003.4138
$ws.Range("L61").Value = 5190.6665
$ws.Range("M61").Value = -2791.4138
$ws.Range("N61").Value = -5614.6665

$ws.Range("H74").Value = 2215.691
$ws.Range("I74").Value = 1252.2157
$ws.Range("J74").Value = 14500
$ws.Range("K74").Value = 1252.2157
$ws.Range("L74").Value = 14500
$ws.Range("M74").Value = -378.2157
$ws.Range("N74").Value = -16248

$ws.Range("H77").Value = 2215.691
$ws.Range("I77").Value = 1252.2157
$ws.Range("J77").Value = 14500
$ws.Range("K77").Value = 6261.0785
$ws.Range("L77").Value = 72500
$ws.Range("M77").Value = -1893.0785
$ws.Range("N77").Value = -81236

$ws.Range("H122").Value = 336324.06
$ws.Range("I122").Value = 2186.8928
$ws.Range("J122").Value = 1004598.44
$ws.Range("K122").Value = 6560.678400000001
$ws.Range("L122").Value = 3013795.32
$ws.Range("M122").Value = -4110.678400000001
$ws.Range("N122").Value = -3018695.32

$ws.Range("H132").Value = 3659.077
$ws.Range("I132").Value = 3112.2
$ws.Range("J132").Value = 4404.8184
$ws.Range("K132").Value = 9336.599999999999
$ws.Range("L132").Value = 13214.4552
$ws.Range("M132").Value = -6806.599999999999
$ws.Range("N132").Value = -18274.4552

$ws.Range("H136").Value = 3643.5854
$ws.Range("I136").Value = 3003.4138
$ws.Range("J136").Value = 5190.6665
$ws.Range("K136").Value = 9010.241399999999
$ws.Range("L136").Value = 15571.9995
$ws.Range("M136").Value = -6460.241399999999
$ws.Range("N136").Value = -20671.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H134").Value = 2672.6765
$ws.Range("I134").Value = 1992.5625
$ws.Range("J134").Value = 4304.95
$ws.Range("K134").Value = 5977.6875
$ws.Range("L134").Value = 12914.85
$ws.Range("M134").Value = -3442.6875
$ws.Range("N134").Value = -17984.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6944.4443
$ws.Range("I31").Value = 5250
$ws.Range("J31").Value = 7428.5713
$ws.Range("K31").Value = 5250
$ws.Range("L31").Value = 7428.5713
$ws.Range("M31").Value = -4955
$ws.Range("N31").Value = -8018.5713

$ws.Range("H34").Value = 6944.4443
$ws.Range("I34").Value = 5250
$ws.Range("J34").Value = 7428.5713
$ws.Range("K34").Value = 5250
$ws.Range("L34").Value = 7428.5713
$ws.Range("M34").Value = -5048
$ws.Range("N34").Value = -7832.5713

$ws.Range("H80").Value = 44933
$ws.Range("J80").Value = 44933
$ws.Range("L80").Value = 44933
$ws.Range("N80").Value = -47179

$ws.Range("H83").Value = 44933
$ws.Range("J83").Value = 44933
$ws.Range("L83").Value = 134799
$ws.Range("N83").Value = -146031

$ws.Range("H86").Value = 11172.963
$ws.Range("I86").Value = 9760.056
$ws.Range("K86").Value = 9760.056
$ws.Range("M86").Value = -8637.056

$ws.Range("H89").Value = 11172.963
$ws.Range("I89").Value = 9760.056
$ws.Range("K89").Value = 48800.28
$ws.Range("M89").Value = -43184.28

$ws.Range("H94").Value = 1945.7646
$ws.Range("J94").Value = 1911.5
$ws.Range("L94").Value = 1911.5
$ws.Range("N94").Value = -2813.5

$ws.Range("H134").Value = 3293591.8
$ws.Range("I134").Value = 3293591.8
$ws.Range("K134").Value = 9880775.399999999
$ws.Range("M134").Value = -9878240.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1471.5769
$ws.Range("J38").Value = 2638.3076
$ws.Range("L38").Value = 7914.9228
$ws.Range("N38").Value = -8608.9228

$ws.Range("H51").Value = 1909.1052
$ws.Range("I51").Value = 802.5
$ws.Range("J51").Value = 3806.1428
$ws.Range("K51").Value = 2407.5
$ws.Range("L51").Value = 11418.4284
$ws.Range("M51").Value = -1947.5
$ws.Range("N51").Value = -12338.4284

$ws.Range("H107").Value = 1151.2106
$ws.Range("I107").Value = 222.5
$ws.Range("J107").Value = 1826.6364
$ws.Range("K107").Value = 667.5
$ws.Range("L107").Value = 5479.9092
$ws.Range("M107").Value = 1252.5
$ws.Range("N107").Value = -9319.9092

$ws.Range("H113").Value = 2291.5833
$ws.Range("I113").Value = 1200.5
$ws.Range("J113").Value = 2837.125
$ws.Range("K113").Value = 3601.5
$ws.Range("L113").Value = 8511.375
$ws.Range("M113").Value = -1431.5
$ws.Range("N113").Value = -12851.375

$ws.Range("H127").Value = 1620
$ws.Range("J127").Value = 1620
$ws.Range("L127").Value = 4860
$ws.Range("N127").Value = -14780

$ws.Range("H129").Value = 3323
$ws.Range("J129").Value = 4546.5557
$ws.Range("L129").Value = 13639.6671
$ws.Range("N129").Value = -23639.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 28091.334
$ws.Range("J15").Value = 28091.334
$ws.Range("L15").Value = 28091.334
$ws.Range("N15").Value = -28667.334

$ws.Range("H70").Value = 5994.5
$ws.Range("I70").Value = 5989
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 5989
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -5719
$ws.Range("N70").Value = -6540

$ws.Range("H73").Value = 5994.5
$ws.Range("I73").Value = 5989
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 5989
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -5053
$ws.Range("N73").Value = -7872

$ws.Range("H81").Value = 28091.334
$ws.Range("J81").Value = 28091.334
$ws.Range("L81").Value = 28091.334
$ws.Range("N81").Value = -30087.334

$ws.Range("H84").Value = 28091.334
$ws.Range("J84").Value = 28091.334
$ws.Range("L84").Value = 84274.00199999999
$ws.Range("N84").Value = -94258.00199999999

$ws.Range("H136").Value = 108021.664
$ws.Range("J136").Value = 108021.664
$ws.Range("L136").Value = 324064.992
$ws.Range("N136").Value = -329164.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4191.846
$ws.Range("J46").Value = 7332.5
$ws.Range("L46").Value = 7332.5
$ws.Range("N46").Value = -7708.5

$ws.Range("H55").Value = 745.3333
$ws.Range("J55").Value = 498.5
$ws.Range("L55").Value = 498.5
$ws.Range("N55").Value = -844.5

$ws.Range("H80").Value = 74999.5
$ws.Range("J80").Value = 74999.5
$ws.Range("L80").Value = 74999.5
$ws.Range("N80").Value = -77245.5

$ws.Range("H83").Value = 74999.5
$ws.Range("J83").Value = 74999.5
$ws.Range("L83").Value = 224998.5
$ws.Range("N83").Value = -236230.5

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 6173.5625
$ws.Range("I136").Value = 2348.375
$ws.Range("J136").Value = 9998.75
$ws.Range("K136").Value = 7045.125
$ws.Range("L136").Value = 29996.25
$ws.Range("M136").Value = -4495.125
$ws.Range("N136").Value = -35096.25

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10838.868
$ws.Range("I132").Value = 13163
$ws.Range("K132").Value = 39489
$ws.Range("M132").Value = -36959
